$wb = $excel.ActiveWorkbook

# The update applies to the 1er Parcial and 3er Parcial sheets only
# (2o Parcial stays as-is). Rows 8 and 9 get updated grade statistics
# now that exams for those two groups/subjects have been scored.
$sheetNames = @("1er Parcial", "3er Parcial")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 8: Polanco Domínguez Rosa María - 3ARHV
    $ws.Range("E8").Value = 25
    $ws.Range("F8").Value = 10
    $ws.Range("G8").Value = 71.43000000000001
    $ws.Range("H8").Value = 28.57
    $ws.Range("I8").Value = 10
    $ws.Range("J8").Value = 10
    $ws.Range("K8").Value = 28.57

    # Row 9: Molina Quezada Raúl - 5ARHV
    $ws.Range("E9").Value = 4
    $ws.Range("F9").Value = 33
    $ws.Range("G9").Value = 10.81
    $ws.Range("H9").Value = 89.19
    $ws.Range("I9").Value = 10
    $ws.Range("J9").Value = 33
    $ws.Range("K9").Value = 89.19
}
